$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117, shifting existing rows 117-141 down to 118-142
$ws.Range("A117").EntireRow.Insert()

# Populate the newly inserted row 117 with the new record's data
$ws.Range("A117").Value = 3
$ws.Range("B117").Value = "Femacal de La Calera"
$ws.Range("C117").Value = "Coquimbo"
$ws.Range("D117").Value = 44722
$ws.Range("E117").Value = 5
$ws.Range("F117").Value = 100112026
$ws.Range("G117").Value = "Haba"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 85
$ws.Range("K117").Value = 22000
$ws.Range("L117").Value = 23000
$ws.Range("M117").Value = 22471
$ws.Range("N117").Value = '$/saco 25 kilos'
$ws.Range("O117").Value = 'Provincia de Limarí'
$ws.Range("P117").Value = 899
$ws.Range("Q117").Value = 25
$ws.Range("R117").Value = "Hortaliza"
